$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Arveja Verde" at Vega Modelo de
# Temuco. It belongs at the top of the existing data block (row 80, just
# after the header + prior entries), so insert a blank row there and push
# everything else down by one.
$ws.Rows(80).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A80").Value2 = 10
$ws.Range("B80").Value2 = "Vega Modelo de Temuco"
$ws.Range("C80").Value2 = "La Araucanía"
$ws.Range("D80").Value2 = 45006
$ws.Range("E80").Value2 = 9
$ws.Range("F80").Value2 = 100112022
$ws.Range("G80").Value2 = "Arveja Verde"
$ws.Range("H80").Value2 = "Sin especificar"
$ws.Range("I80").Value2 = "Primera"
$ws.Range("J80").Value2 = 25
$ws.Range("K80").Value2 = 33000
$ws.Range("L80").Value2 = 33000
$ws.Range("M80").Value2 = 33000
$ws.Range("N80").Value2 = "$/saco 25 kilos"
$ws.Range("O80").Value2 = "Región de La Araucanía"
$ws.Range("P80").Value2 = 1320
$ws.Range("Q80").Value2 = 25
$ws.Range("R80").Value2 = "Hortaliza"
